$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing the existing row 4 (and below) down.
# This turns the old row 4 into row 5, preserving all of its values/formatting.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with this week's data (2021-09-29 / serial 44468).
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44468
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 29000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 29500
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 2950
$ws.Range("T4").Value = 10

Write-Output "done"
